{"js": "// Update the division-practice table: 25 \"dividend\u00f7divisor=quotient, remainder\"\n// cells (the 5 content rows of the 20-row/5-col grid) get new problems.\n// Cells are addressed by (row, col) rather than by searching for the old\n// text, because several old/new strings repeat across different cells\n// (e.g. \"77\u00f73=25, 2\" is both a new value in one cell and an old value in\n// another), which would make a plain text search-and-replace ambiguous.\nconst updates = [\n  { row: 0, col: 0, oldText: \"67\u00f77=9, 4\", newText: \"60\u00f79=6, 6\" },\n  { row: 0, col: 1, oldText: \"75\u00f77=10, 5\", newText: \"10\u00f74=2, 2\" },\n  { row: 0, col: 2, oldText: \"12\u00f73=4, 0\", newText: \"77\u00f73=25, 2\" },\n  { row: 0, col: 3, oldText: \"38\u00f76=6, 2\", newText: \"62\u00f76=10, 2\" },\n  { row: 0, col: 4, oldText: \"19\u00f76=3, 1\", newText: \"69\u00f73=23, 0\" },\n  { row: 4, col: 0, oldText: \"75\u00f75=15, 0\", newText: \"68\u00f78=8, 4\" },\n  { row: 4, col: 1, oldText: \"71\u00f79=7, 8\", newText: \"51\u00f75=10, 1\" },\n  { row: 4, col: 2, oldText: \"81\u00f75=16, 1\", newText: \"85\u00f76=14, 1\" },\n  { row: 4, col: 3, oldText: \"19\u00f74=4, 3\", newText: \"29\u00f73=9, 2\" },\n  { row: 4, col: 4, oldText: \"40\u00f77=5, 5\", newText: \"28\u00f75=5, 3\" },\n  { row: 8, col: 0, oldText: \"52\u00f77=7, 3\", newText: \"42\u00f78=5, 2\" },\n  { row: 8, col: 1, oldText: \"41\u00f77=5, 6\", newText: \"92\u00f74=23, 0\" },\n  { row: 8, col: 2, oldText: \"57\u00f75=11, 2\", newText: \"36\u00f79=4, 0\" },\n  { row: 8, col: 3, oldText: \"77\u00f73=25, 2\", newText: \"64\u00f79=7, 1\" },\n  { row: 8, col: 4, oldText: \"71\u00f74=17, 3\", newText: \"95\u00f75=19, 0\" },\n  { row: 12, col: 0, oldText: \"78\u00f73=26, 0\", newText: \"89\u00f77=12, 5\" },\n  { row: 12, col: 1, oldText: \"74\u00f78=9, 2\", newText: \"55\u00f79=6, 1\" },\n  { row: 12, col: 2, oldText: \"94\u00f74=23, 2\", newText: \"46\u00f78=5, 6\" },\n  { row: 12, col: 3, oldText: \"10\u00f74=2, 2\", newText: \"56\u00f73=18, 2\" },\n  { row: 12, col: 4, oldText: \"27\u00f76=4, 3\", newText: \"68\u00f75=13, 3\" },\n  { row: 16, col: 0, oldText: \"88\u00f72=44, 0\", newText: \"50\u00f74=12, 2\" },\n  { row: 16, col: 1, oldText: \"16\u00f76=2, 4\", newText: \"52\u00f78=6, 4\" },\n  { row: 16, col: 2, oldText: \"85\u00f74=21, 1\", newText: \"32\u00f76=5, 2\" },\n  { row: 16, col: 3, oldText: \"97\u00f72=48, 1\", newText: \"74\u00f78=9, 2\" },\n  { row: 16, col: 4, oldText: \"46\u00f72=23, 0\", newText: \"93\u00f72=46, 1\" },\n];\n\nconst table = context.document.body.tables.getFirstOrNullObject();\nawait context.sync();\n\nconst cells = updates.map((u) => table.getCell(u.row, u.col));\ncells.forEach((c) => c.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < updates.length; i++) {\n  const u = updates[i];\n  const cell = cells[i];\n  const current = (cell.value || \"\").trim();\n  // Guard: only rewrite if the cell still holds the text we expect;\n  // otherwise leave it untouched rather than risk corrupting content.\n  if (current === u.oldText) {\n    cell.value = u.newText;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-practice table: 25 \"dividend\u00f7divisor=quotient, remainder\"\n# cells (the 5 content rows of the 20-row/5-col grid) get new problems.\n# Cells are addressed by (row, col) rather than by searching for the old\n# text, because several old/new strings repeat across different cells\n# (e.g. \"77\u00f73=25, 2\" is both a new value in one cell and an old value in\n# another), which would make a plain text search-and-replace ambiguous.\n$updates = @(\n    @{ Row = 1; Col = 1; OldText = \"67\u00f77=9, 4\"; NewText = \"60\u00f79=6, 6\" },\n    @{ Row = 1; Col = 2; OldText = \"75\u00f77=10, 5\"; NewText = \"10\u00f74=2, 2\" },\n    @{ Row = 1; Col = 3; OldText = \"12\u00f73=4, 0\"; NewText = \"77\u00f73=25, 2\" },\n    @{ Row = 1; Col = 4; OldText = \"38\u00f76=6, 2\"; NewText = \"62\u00f76=10, 2\" },\n    @{ Row = 1; Col = 5; OldText = \"19\u00f76=3, 1\"; NewText = \"69\u00f73=23, 0\" },\n    @{ Row = 5; Col = 1; OldText = \"75\u00f75=15, 0\"; NewText = \"68\u00f78=8, 4\" },\n    @{ Row = 5; Col = 2; OldText = \"71\u00f79=7, 8\"; NewText = \"51\u00f75=10, 1\" },\n    @{ Row = 5; Col = 3; OldText = \"81\u00f75=16, 1\"; NewText = \"85\u00f76=14, 1\" },\n    @{ Row = 5; Col = 4; OldText = \"19\u00f74=4, 3\"; NewText = \"29\u00f73=9, 2\" },\n    @{ Row = 5; Col = 5; OldText = \"40\u00f77=5, 5\"; NewText = \"28\u00f75=5, 3\" },\n    @{ Row = 9; Col = 1; OldText = \"52\u00f77=7, 3\"; NewText = \"42\u00f78=5, 2\" },\n    @{ Row = 9; Col = 2; OldText = \"41\u00f77=5, 6\"; NewText = \"92\u00f74=23, 0\" },\n    @{ Row = 9; Col = 3; OldText = \"57\u00f75=11, 2\"; NewText = \"36\u00f79=4, 0\" },\n    @{ Row = 9; Col = 4; OldText = \"77\u00f73=25, 2\"; NewText = \"64\u00f79=7, 1\" },\n    @{ Row = 9; Col = 5; OldText = \"71\u00f74=17, 3\"; NewText = \"95\u00f75=19, 0\" },\n    @{ Row = 13; Col = 1; OldText = \"78\u00f73=26, 0\"; NewText = \"89\u00f77=12, 5\" },\n    @{ Row = 13; Col = 2; OldText = \"74\u00f78=9, 2\"; NewText = \"55\u00f79=6, 1\" },\n    @{ Row = 13; Col = 3; OldText = \"94\u00f74=23, 2\"; NewText = \"46\u00f78=5, 6\" },\n    @{ Row = 13; Col = 4; OldText = \"10\u00f74=2, 2\"; NewText = \"56\u00f73=18, 2\" },\n    @{ Row = 13; Col = 5; OldText = \"27\u00f76=4, 3\"; NewText = \"68\u00f75=13, 3\" },\n    @{ Row = 17; Col = 1; OldText = \"88\u00f72=44, 0\"; NewText = \"50\u00f74=12, 2\" },\n    @{ Row = 17; Col = 2; OldText = \"16\u00f76=2, 4\"; NewText = \"52\u00f78=6, 4\" },\n    @{ Row = 17; Col = 3; OldText = \"85\u00f74=21, 1\"; NewText = \"32\u00f76=5, 2\" },\n    @{ Row = 17; Col = 4; OldText = \"97\u00f72=48, 1\"; NewText = \"74\u00f78=9, 2\" },\n    @{ Row = 17; Col = 5; OldText = \"46\u00f72=23, 0\"; NewText = \"93\u00f72=46, 1\" }\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nforeach ($u in $updates) {\n    $cell = $tbl.Cell($u.Row, $u.Col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    # Guard: only rewrite if the cell still holds the text we expect;\n    # otherwise leave it untouched rather than risk corrupting content.\n    if ($current -eq $u.OldText) {\n        $cell.Range.Text = $u.NewText\n    }\n}\n"}
